# "added depth into the recursion"
# Re-run of the recursion-depth test harness: the result table now only
# reports codes (4,2)-(6,2) (columns F:H for codes (8,2)-(10,2) are gone),
# the (4,2) column itself was cleared of its header/data text, and the
# remaining columns were refreshed with the new recursion results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook was protected before; it no longer is.
$wb.Unprotect()

# Columns F:H ("Код (8,2)" .. "Код (10,2)") are removed entirely.
$ws.Range("F1:H12").Delete()

# Column B ("Код (4,2)") keeps its header cell (and formatting) but the
# text/values in it are gone.
$ws.Range("B1:B12").ClearContents()

# Rows 8:12 no longer have a value in column D (that recursion branch no
# longer reaches that depth).
$ws.Range("D8:D12").ClearContents()

# Refresh the header row for the remaining result columns.
$ws.Range("C1").Value = "Код (4,2)"
$ws.Range("D1").Value = "Код (5,2)"
$ws.Range("E1").Value = "Код (6,2)"

# New recursion results.
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 31
$ws.Range("E3").Value = 84

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 59

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 33

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 18

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 7

$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 5

$ws.Range("C9").Value = 7
$ws.Range("E9").Value = 1

$ws.Range("C10").Value = 100
$ws.Range("E10").Value = 0

$ws.Range("E11").Value = 0

$ws.Range("E12").Value = 0

# Leave the selection where the author left it when they saved the file.
$ws.Range("I5").Select()
